# Economic Dashboard weekly data refresh - 2025-12-31
# Updates dates / values of several FRED-sourced series and shifts the
# "last updated" highlight (yellow fill, style 49) to the newest date cell
# while restoring the previously-highlighted cell back to the plain
# date style (style 48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: move the "recently updated" yellow-highlight style from one
# date cell to another by copying formats from existing cells that
# already carry style 48 (plain) / 49 (yellow highlight). This re-uses
# the existing style indices instead of fabricating new ones.
# ---------------------------------------------------------------------
$plainDateSource = $ws.Range("C13")   # known style 48 cell (no fill)
$highlightDateSource = $ws.Range("N13")  # known style 49 cell (yellow fill)

# C7: remove highlight (49 -> 48)
$plainDateSource.Copy()
$ws.Range("C7").PasteSpecial(-4122)

# C8: add highlight (48 -> 49) and set new date
$highlightDateSource.Copy()
$ws.Range("C8").PasteSpecial(-4122)

# N51: remove highlight (49 -> 48)
$plainDateSource.Copy()
$ws.Range("N51").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Row 8 - Smoothed Recession Prob (RECPROUSM156N)
# ---------------------------------------------------------------------
$ws.Range("C8").Value2 = 45962
$ws.Range("F8").Value2 = 0.9399999999999999
$ws.Range("G8").Value2 = 0.9
$ws.Range("H8").Value2 = 0.72
$ws.Range("I8").Value2 = 0.7
$ws.Range("J8").Value2 = 0.32

# ---------------------------------------------------------------------
# Row 13 - UI Initial Claims (ICSA)
# ---------------------------------------------------------------------
$ws.Range("N13").Value2 = 46013
$ws.Range("Q13").Value2 = 199000
$ws.Range("R13").Value2 = 215000
$ws.Range("S13").Value2 = 224000
$ws.Range("T13").Value2 = 237000
$ws.Range("U13").Value2 = 192000

# ---------------------------------------------------------------------
# Row 14 - UI Continuing Claims (CCSA)
# ---------------------------------------------------------------------
$ws.Range("N14").Value2 = 46006
$ws.Range("Q14").Value2 = 1866000
$ws.Range("R14").Value2 = 1913000
$ws.Range("S14").Value2 = 1885000
$ws.Range("T14").Value2 = 1830000
$ws.Range("U14").Value2 = 1937000

# ---------------------------------------------------------------------
# Row 29 - 5yr, 5yr Forward (T5YIFR)
# ---------------------------------------------------------------------
$ws.Range("N29").Value2 = 46021
$ws.Range("Q29").Value2 = 2.23
$ws.Range("R29").Value2 = 2.21
$ws.Range("U29").Value2 = 2.24

# ---------------------------------------------------------------------
# Row 30 - 10yr TIPS (T10YIE)
# ---------------------------------------------------------------------
$ws.Range("N30").Value2 = 46021
$ws.Range("Q30").Value2 = 2.24
$ws.Range("R30").Value2 = 2.22
$ws.Range("S30").Value2 = 2.23
$ws.Range("U30").Value2 = 2.24

# ---------------------------------------------------------------------
# Row 47 - FFR (DFF)
# ---------------------------------------------------------------------
$ws.Range("N47").Value2 = 46020

# ---------------------------------------------------------------------
# Row 48 - 2y UST (DGS2)
# ---------------------------------------------------------------------
$ws.Range("N48").Value2 = 46020
$ws.Range("Q48").Value2 = 3.45
$ws.Range("R48").Value2 = 3.46
$ws.Range("S48").Value2 = 3.47
$ws.Range("T48").Value2 = 3.48
$ws.Range("U48").Value2 = 3.44

# ---------------------------------------------------------------------
# Row 49 - 5y UST (DGS5)
# ---------------------------------------------------------------------
$ws.Range("N49").Value2 = 46020
$ws.Range("Q49").Value2 = 3.67
$ws.Range("R49").Value2 = 3.68
$ws.Range("S49").Value2 = 3.7
$ws.Range("T49").Value2 = 3.72
$ws.Range("U49").Value2 = 3.71

# ---------------------------------------------------------------------
# Row 50 - 10y UST (DGS10)
# ---------------------------------------------------------------------
$ws.Range("N50").Value2 = 46020
$ws.Range("Q50").Value2 = 4.12
$ws.Range("R50").Value2 = 4.14
$ws.Range("S50").Value2 = 4.15
$ws.Range("T50").Value2 = 4.18
$ws.Range("U50").Value2 = 4.17

# ---------------------------------------------------------------------
# Row 52 - BAA (DBAA)
# ---------------------------------------------------------------------
$ws.Range("N52").Value2 = 46020
$ws.Range("Q52").Value2 = 5.88
$ws.Range("R52").Value2 = 5.89
$ws.Range("S52").Value2 = 5.88
$ws.Range("T52").Value2 = 5.92
$ws.Range("U52").Value2 = 5.93
